# Learning Diary update (3.5)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the word "tomorrow" from the 30.4.2025 entry.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "that’s what I’m going to focus tomorrow. So far so good",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "that’s what I’m going to focus. So far so good", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "1.5.2025" -> "3.5.2025", kept as two runs ("3" + ".5.2025") like
#    the rest of the dates in this document.
# ---------------------------------------------------------------------
$datePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "1.5.2025") {
        $datePara = $cand
        break
    }
}
$full = $datePara.Range
$firstDigit = $d.Range($full.Start, $full.Start + 1)
$firstDigit.Text = "3"
# Forcing a (no-op) direct-formatting toggle right after the text swap
# makes the run boundary real, matching the multi-run pattern already
# used for the other dates in this document.
$firstDigit2 = $d.Range($full.Start, $full.Start + 1)
$firstDigit2.Bold = 1
$firstDigit2.Bold = 0

# ---------------------------------------------------------------------
# 3) Add the 3.5.2025 diary entry into the (until now empty) paragraph
#    that follows the date paragraph.
# ---------------------------------------------------------------------
$datePara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "3.5.2025") {
        $datePara2 = $cand
        break
    }
}

$part1 = "Today I went through Middleware, Cleanup (Middleware and Handlers) and GET Req Body for POST. Took me some time to understand and going back in the video material quite a few times, but at the end I understood. I stopped at "
$part2 = "GET Req Body for POST"
$part3 = ", so tomorrow I’m starting on File System Module onwards."

# Split off a fresh paragraph right after the date line so the new text
# inherits the correct run formatting (en-US language), then fold the
# extra paragraph mark back out so the target (already-existing, empty)
# paragraph ends up holding the new runs.
$dateRng = $datePara2.Range
$dateRng.InsertParagraphAfter()

$newParaIndex = $datePara2.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Range.InsertBefore($part1 + $part2 + $part3)

$newPara2 = $d.Paragraphs($newParaIndex)
$entryStart = $newPara2.Range.Start
$b1 = $entryStart + $part1.Length
$b2 = $b1 + $part2.Length

$seg1 = $d.Range($entryStart, $b1)
$seg1.Bold = 1
$seg1.Bold = 0
$seg2 = $d.Range($b1, $b2)
$seg2.Bold = 1
$seg2.Bold = 0

# Merge the helper paragraph mark away so the text lands inside the
# original empty paragraph (keeping its identity / following empty
# paragraph untouched).
$newPara3 = $d.Paragraphs($newParaIndex)
$mergeRange = $d.Range($newPara3.Range.End - 1, $newPara3.Range.End)
$mergeRange.Delete()

"done"
